# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
#
# Some "Price" strings look numeric (e.g. "6.80", "0.999") even though the
# sheet stores every value in these columns as plain text. Plain
# `Range.Value = "<numeric-looking string>"` would auto-coerce those into
# real numbers (and drop formatting like the trailing zero), so for those
# cells we briefly force a text number-format, assign, then restore the
# cell's original style to avoid leaving formatting changes behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '42.699.66'
$ws.Range('E2').Value = '  -0.46%  '

$ws.Range('D3').Value = '2.298.03'
$ws.Range('E3').Value = '  -0.09%  '

Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  -0.02%  '

Set-TextValue $ws.Range('D5') '301.07'
$ws.Range('E5').Value = '  -1.39%  '

Set-TextValue $ws.Range('D6') '96.02'
$ws.Range('E6').Value = '  -1.37%  '

$ws.Range('E7').Value = '  -0.98%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('E9').Value = '  -1.64%  '

Set-TextValue $ws.Range('D10') '34.75'
$ws.Range('E10').Value = '  -2.61%  '

Set-TextValue $ws.Range('D11') '19.23'
$ws.Range('E11').Value = '  +4.91%  '

$ws.Range('E12').Value = '  -1.31%  '

$ws.Range('E13').Value = '  +0.11%  '

Set-TextValue $ws.Range('D14') '6.80'
$ws.Range('E14').Value = '  +0.31%  '

$ws.Range('D15').Value = '2.647.69'
$ws.Range('E15').Value = '  -0.41%  '

$ws.Range('D16').Value = '2.287.12'
$ws.Range('E16').Value = '  -0.49%  '

$ws.Range('E17').Value = '  -0.16%  '

$ws.Range('D18').Value = '42.610.46'
$ws.Range('E18').Value = '  -0.53%  '

Set-TextValue $ws.Range('D19') '12.26'
$ws.Range('E19').Value = '  -5.99%  '

$ws.Range('D20').Value = '0.0₃0891'
$ws.Range('E20').Value = '  -1.69%  '

$ws.Range('E21').Value = '  -0.61%  '

Set-TextValue $ws.Range('D22') '67.63'
$ws.Range('E22').Value = '  +0.09%  '

Set-TextValue $ws.Range('D23') '235.37'
$ws.Range('E23').Value = '  -0.51%  '

$ws.Range('E24').Value = '  +3.00%  '

$ws.Range('E25').Value = '  +0.06%  '

$ws.Range('E26').Value = '  -2.54%  '

Set-TextValue $ws.Range('D27') '24.63'
$ws.Range('E27').Value = '  -3.42%  '

Set-TextValue $ws.Range('D28') '2.06'
$ws.Range('E28').Value = '  -0.12%  '

Set-TextValue $ws.Range('D29') '164.45'
$ws.Range('E29').Value = '  -1.56%  '

Set-TextValue $ws.Range('D30') '9.04'
$ws.Range('E30').Value = '  -0.51%  '

Set-TextValue $ws.Range('D31') '32.32'
$ws.Range('E31').Value = '  -2.02%  '

Set-TextValue $ws.Range('D32') '0.999'
$ws.Range('E32').Value = '  -0.05%  '

$ws.Range('E33').Value = '  -1.18%  '

Set-TextValue $ws.Range('D34') '17.46'
$ws.Range('E34').Value = '  +0.65%  '

Set-TextValue $ws.Range('D35') '4.44'
$ws.Range('E35').Value = '  -7.66%  '

$ws.Range('E36').Value = '  +0.58%  '

$ws.Range('E37').Value = '  -2.84%  '

$ws.Range('E38').Value = '  -1.91%  '

$ws.Range('E39').Value = '  -0.05%  '

$ws.Range('E40').Value = '  -0.42%  '

$ws.Range('E41').Value = '  -1.27%  '

Set-TextValue $ws.Range('D42') '19.74'
$ws.Range('E42').Value = '  +9.15%  '

$ws.Range('D43').Value = '1.953.73'
$ws.Range('E43').Value = '  -3.11%  '

$ws.Range('E44').Value = '  +4.60%  '

$ws.Range('E45').Value = '  -0.83%  '

$ws.Range('E46').Value = '  -4.10%  '

$ws.Range('E47').Value = '  -1.24%  '

Set-TextValue $ws.Range('D48') '2.90'
$ws.Range('E48').Value = '  +0.00%  '

$ws.Range('E49').Value = '  -0.20%  '

Set-TextValue $ws.Range('D50') '53.24'
$ws.Range('E50').Value = '  -1.43%  '

Set-TextValue $ws.Range('D51') '71.57'
$ws.Range('E51').Value = '  -0.77%  '
